$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay stored as text,
# matching the source data which uses literal strings like "1.00" / "5.30".
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.184.04'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '3.464.52'
$ws.Range("E3").Value = '  +2.27%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '578.74'
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").Value = '148.41'
$ws.Range("E6").Value = '  +7.12%  '

$ws.Range("D7").Value = '3.465.47'
$ws.Range("E7").Value = '  +2.36%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '0.476'
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").Value = '7.71'
$ws.Range("E10").Value = '  +2.93%  '

$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +0.09%  '

$ws.Range("E12").Value = '  -0.51%  '

$ws.Range("D13").Value = '4.055.18'
$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").Value = '27.50'
$ws.Range("E15").Value = '  +5.34%  '

$ws.Range("D16").Value = '0.0000177'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").Value = '3.461.07'
$ws.Range("E17").Value = '  +2.36%  '

$ws.Range("D18").Value = '62.100.46'
$ws.Range("E18").Value = '  +0.39%  '

$ws.Range("D19").Value = '6.22'
$ws.Range("E19").Value = '  +4.78%  '

$ws.Range("D20").Value = '14.19'
$ws.Range("E20").Value = '  +1.21%  '

$ws.Range("D21").Value = '9.59'
$ws.Range("E21").Value = '  +1.60%  '

$ws.Range("D22").Value = '391.23'
$ws.Range("E22").Value = '  +3.27%  '

$ws.Range("D23").Value = '0.565'
$ws.Range("E23").Value = '  +1.26%  '

$ws.Range("D24").Value = '3.587.87'
$ws.Range("E24").Value = '  +2.00%  '

$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = '72.68'
$ws.Range("E26").Value = '  +1.56%  '

$ws.Range("D27").Value = '5.77'
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").Value = '0.0000126'
$ws.Range("E28").Value = '  -0.08%  '

$ws.Range("D29").Value = '0.178'
$ws.Range("E29").Value = '  +7.89%  '

$ws.Range("D30").Value = '7.84'
$ws.Range("E30").Value = '  +2.21%  '

$ws.Range("E31").Value = '  -13.19%  '

$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("D33").Value = '8.29'
$ws.Range("E33").Value = '  -0.37%  '

$ws.Range("D34").Value = '2.18'
$ws.Range("E34").Value = '  +0.48%  '

$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").Value = '24.08'
$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("D37").Value = '5.33'
$ws.Range("E37").Value = '  +1.39%  '

$ws.Range("D38").Value = '7.08'
$ws.Range("E38").Value = '  +3.09%  '

$ws.Range("E39").Value = '  +1.65%  '

$ws.Range("D40").Value = '166.48'
$ws.Range("E40").Value = '  +0.95%  '

$ws.Range("D41").Value = '0.0793'
$ws.Range("E41").Value = '  +2.46%  '

$ws.Range("D42").Value = '26.38'
$ws.Range("E42").Value = '  +7.84%  '

$ws.Range("D43").Value = '0.797'
$ws.Range("E43").Value = '  +2.61%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.75'
$ws.Range("E44").Value = '  +0.37%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '42.28'
$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("D47").Value = '4.48'
$ws.Range("E47").Value = '  +1.75%  '

$ws.Range("D48").Value = '1.22'
$ws.Range("E48").Value = '  -0.18%  '

$ws.Range("D49").Value = '2.650.37'
$ws.Range("E49").Value = '  +10.81%  '

$ws.Range("D50").Value = '23.76'
$ws.Range("E50").Value = '  +2.91%  '

$ws.Range("D51").Value = '6.91'
$ws.Range("E51").Value = '  +0.36%  '

# Remove the stray number-format/style applied above so the cells keep
# matching their original (unstyled) appearance.
$ws.Range("D2:D51").ClearFormats()